{"js": "// Append the \"Aula 3\" heading (bold + underline, numbering level 0) and a\n// following empty sub-bullet (numbering level 1) after the last paragraph\n// of the document body \u2014 mirrors the existing \"Aula 1\"/\"Aula 2\" headings.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst lastParagraph = paragraphs.items[paragraphs.items.length - 1];\n\n// Insert both new paragraphs first, while they still inherit the plain\n// (non-bold, non-underlined) run formatting of the preceding paragraph \u2014\n// then apply bold/underline only to the heading run afterwards, so the\n// second (blank) paragraph stays plain, matching the target document.\nconst headingParagraph = lastParagraph.insertParagraph(\n  \"Aula 3 \u2013 Desenvolvendo e Testando o Corretor:\",\n  \"After\"\n);\nheadingParagraph.listItemOrNullObject.level = 0;\n\nconst subParagraph = headingParagraph.insertParagraph(\" \", \"After\");\nsubParagraph.listItemOrNullObject.level = 1;\n\nheadingParagraph.font.bold = true;\nheadingParagraph.font.boldBidirectional = true;\nheadingParagraph.font.underline = \"Single\";\n\nawait context.sync();\n", "ps1": "# Append the \"Aula 3\" heading (bold + underline, numbering level 0) and a\n# following empty sub-bullet (numbering level 1) after the last paragraph\n# of the document body \u2014 mirrors the existing \"Aula 1\"/\"Aula 2\" headings.\n\n$d = $word.ActiveDocument\n\n# Insert the heading paragraph right after the current last paragraph.\n$lastIndex = $d.Paragraphs.Count\n$lastParagraph = $d.Paragraphs.Item($lastIndex)\n$lastParagraph.Range.InsertParagraphAfter()\n\n# Insert the (initially empty) sub-bullet paragraph right after the heading.\n$headingIndex = $d.Paragraphs.Count\n$headingParagraph = $d.Paragraphs.Item($headingIndex)\n$headingParagraph.Range.InsertParagraphAfter()\n\n$subIndex = $d.Paragraphs.Count\n\n# Fill in the heading paragraph text + list level, both inherited as plain\n# (non-bold, non-underlined) text from \"O que s\u00e3o tipos de palavras.\" above.\n$headingParagraph = $d.Paragraphs.Item($headingIndex)\n$headingParagraph.Range.Text = \"Aula 3 \u2013 Desenvolvendo e Testando o Corretor:\"\n$headingParagraph.Range.ListFormat.ListLevelNumber = 1\n\n# Bold + underline only the visible text run, not the paragraph mark, so the\n# paragraph-mark run properties (w:pPr/w:rPr) stay plain like the source.\n$headingText = $headingParagraph.Range.Duplicate\n[void]$headingText.MoveEnd(1, -1)\n$headingText.Font.Bold = 1\n$headingText.Font.BoldBi = 1\n$headingText.Font.Underline = 1\n\n# Fill in the blank sub-bullet paragraph text + list level.\n$subParagraph = $d.Paragraphs.Item($subIndex)\n$subParagraph.Range.Text = \" \"\n$subParagraph.Range.ListFormat.ListLevelNumber = 2\n"}
